# Generate Report for Handback
# The d6c199f8 file has been handed back (in sync with en-US), so update its
# status / handback timestamp / error detail across all report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns for the d6c199f8 row
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack

# zh-cn detail sheet: update Status, Latest Handback DateTime, and clear the
# stale Error Detail for the d6c199f8 row
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("K3").Value = "2016-08-13 08:57:56"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).AutoFit()

# de-de detail sheet: same update for the d6c199f8 row
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("K3").Value = "2016-08-13 08:58:10"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).AutoFit()
